# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12): fix numbers & give A10/A11/A12 the mtitleStyle
#     that A9 ("Right"/"Wrong"/... header row) already uses.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 24
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 96
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = "96/112"

# --- Per-question answer grid (rows 16-40): column A becomes the "Student Ans"
#     value, filled in wherever the student answered (copy the green
#     correctStyle formatting already used by B10/B11/B12, then set the text).
#     Rows 21, 23, 28 and 40 are left blank (not attempted), matching B's
#     original formatting there.
$ws.Range("B10").Copy()
$ws.Range("A16:A20").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A24:A27").PasteSpecial(-4122)
$ws.Range("A29:A39").PasteSpecial(-4122)

$ws.Range("A16").Value = "Option A"
$ws.Range("A17").Value = "Option D"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A20").Value = "Option B"
$ws.Range("A22").Value = "Option D"
$ws.Range("A24").Value = "Option A"
$ws.Range("A25").Value = "Option A"
$ws.Range("A26").Value = "Option C"
$ws.Range("A27").Value = "Option A"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A34").Value = "Option B"
$ws.Range("A35").Value = "Option D"
$ws.Range("A36").Value = "Option A"
$ws.Range("A37").Value = "Option A"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"

# --- Drop the 2nd and 3rd students' columns (D:E for the question rows, and
#     the whole G:H block) -- only one student's answers are kept now.
$ws.Range("D16:E40").Clear()
$ws.Range("G15:H40").Clear()
